# Updates symbol price/volume/hour data for Tue Jan 17 07:10:21 UTC 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row number, new Price (col D, optional), new Volume(1h) % (col E, optional),
# new Hora (col G, always "7" this run). $null means "leave this cell untouched".
$updates = @(
    @{ Row = 2; D = '298.87'; E = '-1.36%'; G = '7' },
    @{ Row = 3; D = '31.46'; E = '-1.00%'; G = '7' },
    @{ Row = 4; D = '5.093'; E = '-1.38%'; G = '7' },
    @{ Row = 5; D = '0.07874'; E = '0.78%'; G = '7' },
    @{ Row = 6; D = '2.273'; E = '-6.77%'; G = '7' },
    @{ Row = 7; D = '7.801'; E = '-2.12%'; G = '7' },
    @{ Row = 8; D = '3.855'; E = '-0.38%'; G = '7' },
    @{ Row = 9; D = '0.9188'; E = '0.84%'; G = '7' },
    @{ Row = 10; D = '0.1738'; E = '0.33%'; G = '7' },
    @{ Row = 11; D = '0.07580'; E = '3.23%'; G = '7' },
    @{ Row = 12; D = '0.09379'; E = '15.48%'; G = '7' },
    @{ Row = 13; D = '0.03009'; E = '-0.99%'; G = '7' },
    @{ Row = 14; D = $null; E = '0.67%'; G = '7' },
    @{ Row = 15; D = '0.001506'; E = '-1.30%'; G = '7' },
    @{ Row = 16; D = '0.006062'; E = '-2.18%'; G = '7' },
    @{ Row = 17; D = $null; E = '-0.59%'; G = '7' },
    @{ Row = 18; D = '2.245'; E = '0.13%'; G = '7' },
    @{ Row = 19; D = $null; E = '0.88%'; G = '7' },
    @{ Row = 20; D = '0.1310'; E = '-2.06%'; G = '7' },
    @{ Row = 21; D = '3.952'; E = '-15.79%'; G = '7' },
    @{ Row = 22; D = '0.1712'; E = '9.35%'; G = '7' },
    @{ Row = 23; D = '0.04620'; E = '-0.75%'; G = '7' },
    @{ Row = 24; D = '0.001255'; E = '-0.45%'; G = '7' },
    @{ Row = 25; D = '0.004469'; E = '-1.22%'; G = '7' },
    @{ Row = 26; D = '0.0001250'; E = '-7.32%'; G = '7' },
    @{ Row = 27; D = '0.0003399'; E = '24.00%'; G = '7' },
    @{ Row = 28; D = $null; E = $null; G = '7' },
    @{ Row = 29; D = $null; E = $null; G = '7' },
    @{ Row = 30; D = $null; E = $null; G = '7' },
    @{ Row = 31; D = $null; E = $null; G = '7' },
    @{ Row = 32; D = $null; E = $null; G = '7' },
    @{ Row = 33; D = $null; E = $null; G = '7' },
    @{ Row = 34; D = $null; E = $null; G = '7' },
    @{ Row = 35; D = $null; E = $null; G = '7' },
    @{ Row = 36; D = $null; E = $null; G = '7' },
    @{ Row = 37; D = $null; E = $null; G = '7' },
    @{ Row = 38; D = $null; E = $null; G = '7' },
    @{ Row = 39; D = '0.01732'; E = '-3.25%'; G = '7' },
    @{ Row = 40; D = '0.04606'; E = '0.81%'; G = '7' },
    @{ Row = 41; D = '0.007005'; E = '-3.79%'; G = '7' },
    @{ Row = 42; D = '0.1356'; E = '-0.56%'; G = '7' },
    @{ Row = 43; D = '0.002191'; E = '-2.18%'; G = '7' },
    @{ Row = 44; D = '0.01026'; E = '-4.74%'; G = '7' },
    @{ Row = 45; D = '0.00006267'; E = '-3.02%'; G = '7' },
    @{ Row = 46; D = $null; E = '0.08%'; G = '7' },
    @{ Row = 47; D = '0.007982'; E = '-19.33%'; G = '7' },
    @{ Row = 48; D = '0.7452'; E = '-9.18%'; G = '7' },
    @{ Row = 49; D = '0.00002101'; E = '0.08%'; G = '7' },
    @{ Row = 50; D = '0.0002001'; E = '0.08%'; G = '7' },
    @{ Row = 51; D = $null; E = $null; G = '7' }
)

foreach ($u in $updates) {
    # A leading apostrophe forces Excel to store the value as literal text
    # (matching the source file's inline-string cells) instead of coercing
    # numeric- or percent-looking text into a number.
    if ($null -ne $u.D) { $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D }
    if ($null -ne $u.E) { $ws.Cells.Item($u.Row, 5).Value = "'" + $u.E }
    $ws.Cells.Item($u.Row, 7).Value = "'" + $u.G
}
